$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("A8").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"

$ws.Range("A7").Value = "3"
$ws.Range("B7").Value = "Acacia Breeze"
$ws.Range("C7").Value = "T2345678D"
$ws.Range("D7").Value = "1"
$ws.Range("E7").Value = "Are you there?"
$ws.Range("F7").Value = "S5678901G"
$ws.Range("G7").Value = "Manager"
$ws.Range("H7").Value = "2025-04-23T23:06:54.954202"

$ws.Range("A8").Value = "3"
$ws.Range("B8").Value = "Acacia Breeze"
$ws.Range("C8").Value = "T2345678D"
$ws.Range("D8").Value = "2"
$ws.Range("E8").Value = "test"
$ws.Range("F8").Value = "S5678901G"
$ws.Range("G8").Value = "Manager"
$ws.Range("H8").Value = "2025-04-23T23:07:20.782654"

$ws.Range("A7").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("A8").ClearFormats()
$ws.Range("D8").ClearFormats()
